$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force Text format on the cells we will update so Excel keeps the
# values as literal strings (matching percent/decimal formatting)
# instead of re-interpreting them as numbers.
$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","E16","D17","E17","D19","E19","D20","E20","D21","E21","E22","D23","E23","D24","E24","D25","E25","E26","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","E49","D50","E50","D51","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "302.40"
$ws.Range("E2").Value = "2.10%"
$ws.Range("D3").Value = "44.16"
$ws.Range("E3").Value = "6.63%"
$ws.Range("D4").Value = "5.091"
$ws.Range("E4").Value = "0.96%"
$ws.Range("D5").Value = "0.07705"
$ws.Range("E5").Value = "3.20%"
$ws.Range("D6").Value = "4.418"
$ws.Range("E6").Value = "1.41%"
$ws.Range("E7").Value = "2.51%"
$ws.Range("D8").Value = "1.046"
$ws.Range("E8").Value = "12.76%"
$ws.Range("D9").Value = "0.1276"
$ws.Range("E9").Value = "7.69%"
$ws.Range("D10").Value = "0.1852"
$ws.Range("E10").Value = "1.30%"
$ws.Range("D11").Value = "0.09239"
$ws.Range("E11").Value = "4.15%"
$ws.Range("D12").Value = "0.04173"
$ws.Range("E12").Value = "-0.46%"
$ws.Range("D13").Value = "0.1046"
$ws.Range("E13").Value = "-0.55%"
$ws.Range("D14").Value = "0.001279"
$ws.Range("E14").Value = "0.19%"
$ws.Range("D15").Value = "0.005759"
$ws.Range("E15").Value = "-3.63%"
$ws.Range("E16").Value = "1,911.15%"
$ws.Range("D17").Value = "3.345"
$ws.Range("E17").Value = "0.01%"
$ws.Range("D19").Value = "0.3342"
$ws.Range("E19").Value = "0.99%"
$ws.Range("D20").Value = "8.089"
$ws.Range("E20").Value = "2.58%"
$ws.Range("D21").Value = "0.1360"
$ws.Range("E21").Value = "-3.41%"
$ws.Range("E22").Value = "7.09%"
$ws.Range("D23").Value = "0.04190"
$ws.Range("E23").Value = "3.92%"
$ws.Range("D24").Value = "0.001283"
$ws.Range("E24").Value = "1.49%"
$ws.Range("D25").Value = "0.004419"
$ws.Range("E25").Value = "14.05%"
$ws.Range("E26").Value = "9.60%"
$ws.Range("D38").Value = "0.02495"
$ws.Range("E38").Value = "4.34%"
$ws.Range("D39").Value = "0.05302"
$ws.Range("E39").Value = "1.76%"
$ws.Range("D40").Value = "0.005925"
$ws.Range("E40").Value = "-11.37%"
$ws.Range("D41").Value = "0.007723"
$ws.Range("E41").Value = "-0.92%"
$ws.Range("D42").Value = "0.1351"
$ws.Range("E42").Value = "2.20%"
$ws.Range("D43").Value = "0.007362"
$ws.Range("E43").Value = "-0.11%"
$ws.Range("D44").Value = "0.007536"
$ws.Range("E44").Value = "4.78%"
$ws.Range("D45").Value = "0.3012"
$ws.Range("E45").Value = "-6.29%"
$ws.Range("D46").Value = "0.00006658"
$ws.Range("E46").Value = "6.90%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("D48").Value = "0.04309"
$ws.Range("E48").Value = "-6.36%"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "-0.14%"
